$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '30.554.69'
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").Value = '1.918.89'
$ws.Range("E3").Value = '  -0.06%  '

$ws.Range("E4").Value = '  +0.05%  '

Set-TextValue $ws.Range("D5") '244.80'
$ws.Range("E5").Value = '  +1.65%  '

Set-TextValue $ws.Range("D6") '1.000'
$ws.Range("E6").Value = '  -0.02%  '

Set-TextValue $ws.Range("D7") '0.4877'
$ws.Range("E7").Value = '  +3.85%  '

Set-TextValue $ws.Range("D8") '0.2917'
$ws.Range("E8").Value = '  +1.95%  '

Set-TextValue $ws.Range("D9") '0.06737'
$ws.Range("E9").Value = '  -2.22%  '

Set-TextValue $ws.Range("D10") '107.42'
$ws.Range("E10").Value = '  +0.18%  '

Set-TextValue $ws.Range("D11") '18.73'
$ws.Range("E11").Value = '  +2.02%  '

$ws.Range("D12").Value = '1.927.74'
$ws.Range("E12").Value = '  +0.33%  '

Set-TextValue $ws.Range("D13") '0.07657'
$ws.Range("E13").Value = '  +0.05%  '

Set-TextValue $ws.Range("D14") '5.305'
$ws.Range("E14").Value = '  +2.21%  '

Set-TextValue $ws.Range("D15") '0.6692'
$ws.Range("E15").Value = '  +1.85%  '

Set-TextValue $ws.Range("D16") '278.57'
$ws.Range("E16").Value = '  -5.65%  '

$ws.Range("D17").Value = '30.545.29'
$ws.Range("E17").Value = '  +0.37%  '

$ws.Range("E18").Value = '  +0.09%  '

Set-TextValue $ws.Range("D19") '0.000007566'
$ws.Range("E19").Value = '  -1.10%  '

$ws.Range("D20").Value = '2.168.91'
$ws.Range("E20").Value = '  +0.52%  '

Set-TextValue $ws.Range("D21") '12.82'
$ws.Range("E21").Value = '  -1.68%  '

Set-TextValue $ws.Range("D22") '5.507'
$ws.Range("E22").Value = '  +5.20%  '

Set-TextValue $ws.Range("D23") '1.000'
$ws.Range("E23").Value = '  -0.08%  '

Set-TextValue $ws.Range("D24") '6.430'
$ws.Range("E24").Value = '  +3.46%  '

Set-TextValue $ws.Range("D25") '9.452'
$ws.Range("E25").Value = '  +1.95%  '

Set-TextValue $ws.Range("D26") '164.60'
$ws.Range("E26").Value = '  -1.91%  '

Set-TextValue $ws.Range("D27") '20.23'
$ws.Range("E27").Value = '  -5.86%  '

Set-TextValue $ws.Range("D28") '2.106'
$ws.Range("E28").Value = '  +2.95%  '

Set-TextValue $ws.Range("D29") '0.1055'
$ws.Range("E29").Value = '  -2.73%  '

Set-TextValue $ws.Range("D30") '1.405'
$ws.Range("E30").Value = '  +3.32%  '

Set-TextValue $ws.Range("D31") '4.159'
$ws.Range("E31").Value = '  +0.22%  '

$ws.Range("E32").Value = '  +2.17%  '

Set-TextValue $ws.Range("D33") '0.05004'
$ws.Range("E33").Value = '  -1.49%  '

Set-TextValue $ws.Range("D34") '0.7330'
$ws.Range("E34").Value = '  -1.40%  '

Set-TextValue $ws.Range("D35") '1.139'
$ws.Range("E35").Value = '  -0.57%  '

Set-TextValue $ws.Range("D36") '1.000'

Set-TextValue $ws.Range("D37") '2.729'
$ws.Range("E37").Value = '  -0.47%  '

Set-TextValue $ws.Range("D38") '0.02039'
$ws.Range("E38").Value = '  +1.18%  '

Set-TextValue $ws.Range("D39") '2.678'
$ws.Range("E39").Value = '  -0.28%  '

Set-TextValue $ws.Range("D40") '111.67'
$ws.Range("E40").Value = '  +2.98%  '

$ws.Range("E41").Value = '  -1.92%  '

Set-TextValue $ws.Range("D42") '0.4455'
$ws.Range("E42").Value = '  +5.41%  '

Set-TextValue $ws.Range("D43") '0.8739'
$ws.Range("E43").Value = '  +0.02%  '

Set-TextValue $ws.Range("D44") '5.905'
$ws.Range("E44").Value = '  +1.16%  '

Set-TextValue $ws.Range("D45") '1.000'
$ws.Range("E45").Value = '  -0.01%  '

Set-TextValue $ws.Range("D46") '68.09'
$ws.Range("E46").Value = '  +0.51%  '

$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D47") '7.275'
$ws.Range("E47").Value = '  +0.96%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D48") '9.318'
$ws.Range("E48").Value = '  +0.99%  '

$ws.Range("B49").Value = 'BitcoinSV'
$ws.Range("C49").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue $ws.Range("D49") '48.55'
$ws.Range("E49").Value = '  -8.67%  '

Set-TextValue $ws.Range("D50") '0.1256'
$ws.Range("E50").Value = '  +3.94%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D51") '1.471'
$ws.Range("E51").Value = '  +6.67%  '
